$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.922.99"
$ws.Range("E2").Value = "  +1.40%  "

$ws.Range("D3").Value = "3.810.63"
$ws.Range("E3").Value = "  +0.97%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "445.02"
$ws.Range("E5").Value = "  +6.37%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.81"
$ws.Range("E6").Value = "  +14.85%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.624"
$ws.Range("E7").Value = "  +5.00%  "

$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.738"
$ws.Range("E9").Value = "  +3.65%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.156"
$ws.Range("E10").Value = "  -1.90%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000322"
$ws.Range("E11").Value = "  -5.58%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.53"
$ws.Range("E12").Value = "  +10.53%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.33"
$ws.Range("E13").Value = "  +3.32%  "

$ws.Range("D14").Value = "4.391.22"
$ws.Range("E14").Value = "  +0.71%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.97"
$ws.Range("E15").Value = "  -7.63%  "

$ws.Range("D16").Value = "3.818.24"
$ws.Range("E16").Value = "  +1.31%  "

$ws.Range("E17").Value = "  -0.09%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.97"
$ws.Range("E18").Value = "  +3.85%  "

$ws.Range("E19").Value = "  +8.25%  "

$ws.Range("D20").Value = "66.992.52"
$ws.Range("E20").Value = "  +1.16%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "426.27"
$ws.Range("E21").Value = "  +6.30%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.66"
$ws.Range("E22").Value = "  +4.20%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.24"
$ws.Range("E23").Value = "  +9.96%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "86.71"
$ws.Range("E24").Value = "  +4.97%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "37.44"
$ws.Range("E25").Value = "  +3.00%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.42"
$ws.Range("E26").Value = "  +8.87%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.54"
$ws.Range("E27").Value = "  -2.03%  "

$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.69"
$ws.Range("E28").Value = "  +5.57%  "

$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.45"
$ws.Range("E29").Value = "  +18.78%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "740.48"

$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.67"
$ws.Range("E31").Value = "  +12.82%  "

$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.133"
$ws.Range("E32").Value = "  +12.33%  "

$ws.Range("E33").Value = "  -0.97%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "43.08"
$ws.Range("E34").Value = "  +16.24%  "

$ws.Range("E35").Value = "  +5.90%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "58.06"
$ws.Range("E36").Value = "  +6.70%  "

$ws.Range("E37").Value = "  -0.21%  "

$ws.Range("E38").Value = "  +19.63%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0475"
$ws.Range("E39").Value = "  +6.31%  "

$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.346"
$ws.Range("E40").Value = "  +19.30%  "

$ws.Range("B41").Value = "ThetaToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.88"
$ws.Range("E41").Value = "  -0.77%  "

$ws.Range("E42").Value = "  +0.15%  "

$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.140"
$ws.Range("E43").Value = "  +5.55%  "

$ws.Range("B44").Value = "PEPE"
$ws.Range("C44").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D44").Value = "0.0₃0672"
$ws.Range("E44").Value = "  -10.40%  "

$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.27"
$ws.Range("E45").Value = "  +6.42%  "

$ws.Range("B46").Value = "LidoDAOToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.40"
$ws.Range("E46").Value = "  +4.34%  "

$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.48"
$ws.Range("E47").Value = "  +13.66%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "147.13"
$ws.Range("E48").Value = "  +2.79%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.11"
$ws.Range("E49").Value = "  +5.02%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.66"
$ws.Range("E50").Value = "  +6.90%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.88"
$ws.Range("E51").Value = "  +6.35%  "
